# Delete the row belonging to account 004547722 / MARCIA / 5000
# and shift all rows below it up by one (matching the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = $null

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $acct = $ws.Cells.Item($i, 1).Value2
    $name = $ws.Cells.Item($i, 2).Value2
    if ($acct -eq "004547722" -and $name -eq "MARCIA") {
        $targetRow = $i
        break
    }
}

if ($targetRow -ne $null) {
    $ws.Rows.Item($targetRow).Delete()
}
